# Added Proto-pasta PC-ABS profile that uses cooling fans
#
# This script applies the following changes to the Material profiles sheet:
#  - Renames "Protopasta Conductive" -> "Proto-pasta Conductive" (row 16, col A)
#  - Renames "Protopasta PC-ABS" -> "Proto-pasta PC-ABS" (row 18, col A/B)
#  - Marks the existing Proto-pasta PC-ABS print-settings profile as the
#    fan-cooled variant: "Proto-pasta PC-ABS" -> "Proto-pasta PC-ABS / (fans)" (row 18, col C)
#  - Normalises row 18's printer column to the standard purgebubble profile
#  - Marks the Polymaker PC-Max print-settings profile as the no-fan variant:
#    "Polymaker PC-Max fans" -> "Polymaker PC-Max fans / (no fans)" (row 14, col C)
#  - Adds an explanatory note row under the header, in C2
#  - Updates the D18 cell comment text
#  - Updates the saved selection to E13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Protopasta ..." entries to the hyphenated "Proto-pasta ..." spelling
$ws.Range("A16").Value2 = "Proto-pasta Conductive"
$ws.Range("A18").Value2 = "Proto-pasta PC-ABS"
$ws.Range("B18").Value2 = "Proto-pasta PC-ABS"

# Distinguish the alternate (fan / no-fan) print-settings profiles
$ws.Range("C18").Value2 = "Proto-pasta PC-ABS / (fans)"
$ws.Range("C14").Value2 = "Polymaker PC-Max fans / (no fans)"

# Row 18 now uses the normal purgebubble printer profile
$ws.Range("D18").Value2 = "Original Prusa i3 MK3 purgebubble"

# Add the new explanatory note under the header row
$ws.Range("C2").Value2 = "Alternate profiles denoted with parentheses"

# Update the comment on D18 to reflect that the purgebubble method now works well
$comment = $ws.Range("D18").Comment
[void]$comment.Text("Darragh Broadbent:" + [char]10 + "Brim suggested, prints well with the purgebubble method assuming glue stick used on bed.")

# Restore the saved cell selection
[void]$ws.Range("E13").Select()
